# Apply edits described by the diff:
# - Rename "Jakay N" -> "Jakay M" (Test Data sheet, shared string reused)
# - Add a new "close" test step row to the Test Steps sheet
# - Mark the "verifyText" row (row 7) in Test Steps as "Executed"
# - Make "Test Steps" sheet the active sheet/tab
# - Update selection on Test Steps and Test Data sheets

$wb = $excel.ActiveWorkbook

$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsData  = $wb.Worksheets.Item("Test Data")

# Rename the value in Test Data!D3 from "Jakay N" to "Jakay M"
$wsData.Range("D3").Value = "Jakay M"

# Mark row 7 (verifyText step) as Executed
$wsSteps.Range("F7").Value = "Executed"

# Add a new row 8 for a "close" step - clone the formatting of row 7 first,
# then fill in the values for the new step.
$wsSteps.Range("A7:F7").Copy()
$wsSteps.Range("A8:F8").PasteSpecial(-4122)
$wsSteps.Range("A8").Value = "TC_001_Validate_Login_Page"
$wsSteps.Range("B8").Value = ""
$wsSteps.Range("C8").Value = "close"
$wsSteps.Range("D8").Value = ""
$wsSteps.Range("E8").Value = ""
$wsSteps.Range("F8").Value = "Executed"

# Update selections
$wsData.Range("D3").Select()
$wsSteps.Range("A10").Select()

# Make "Test Steps" the active sheet (so it is saved as the active tab)
$wsSteps.Activate()
